$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the "Date" metadata value.
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"

# Insert a new row for "Jurisdiction" right after "Contact" (row 10),
# pushing Description/Purpose/.../Count down by one row.
$ws.Rows.Item(11).Insert()

# Copy the format of the row above (Contact) onto the newly inserted row
# so the new cells keep the same style (border/alignment) as the rest of
# the table instead of Excel's default insert style.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
